$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 holds three "type: blog" widgets whose trailing "ser: N" post id is
# bumped by 2 (a new post, #165, was published; the chain shifts up).
$ws.Range("H10").Value = "type: blog`nwidth: 2`nheight: 1`nser: 163"
$ws.Range("D10").Value = "type: blog`nwidth: 2`nheight: 1`nser: 164"
$ws.Range("B10").Value = "type: blog`nwidth: 2`nheight: 1`nser: 165"
